# No34 reviewed and No384. Shuffle an Array finished
# Append two new rows (49, 50) to the LeetCode tracking sheet, matching the
# formatting of similar existing rows, then wire up hyperlinks for the new
# "link" column cells and move the active selection like the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 49: 384. Shuffle an Array (new entry, not yet reviewed) ---------
# Style pattern (A..G) matches existing row 33 (s=34,18,19,20,21,21,44);
# column H (wrong-answer mark) matches row 47's style (s=47).
$ws.Range("A33:G33").Copy()
$ws.Range("A49:G49").PasteSpecial(-4122)
$ws.Range("H47").Copy()
$ws.Range("H49").PasteSpecial(-4122)

$ws.Range("A49").Value = "384. Shuffle an Array"
$ws.Range("B49").Value = "Medium"
$ws.Range("C49").Value = "https://leetcode.com/problems/shuffle-an-array/"
$ws.Range("D49").Value = 44554
$ws.Range("E49").Value = "数学"
$ws.Range("F49").Value = "注意本题考查洗牌算法，要求完全随机"
$ws.Range("G49").Value = "未复习"
$ws.Range("H49").Value = "⭕"

# --- Row 50: 34. Find First and Last Position of Element in Sorted Array -
# Style pattern (A..G) matches existing row 32 / row 8 (s=27,8,9,10,11,11,28).
$ws.Range("A32:G32").Copy()
$ws.Range("A50:G50").PasteSpecial(-4122)

$ws.Range("A50").Value = "34. Find First and Last Position of Element in Sorted Array"
$ws.Range("B50").Value = "Medium"
$ws.Range("C50").Value = "https://leetcode.com/problems/find-first-and-last-position-of-element-in-sorted-array/"
$ws.Range("D50").Value = 44468
$ws.Range("E50").Value = "二分法"
$ws.Range("F50").Value = "二分法找元素当存在多个元素时如何定位最小index和最大index；"
$ws.Range("G50").Value = 44554

# --- Hyperlinks for the new "link" column cells ---------------------------
# Hyperlinks.Add() re-styles the anchor cell with a freshly synthesised
# "hyperlink" xf instead of reusing the workbook's existing matching style,
# so re-apply the intended formatting (copied from the same donor rows
# above) on top once the hyperlink relationship is wired up.
$ws.Hyperlinks.Add($ws.Range("C49"), "https://leetcode.com/problems/shuffle-an-array/")
$ws.Range("C33").Copy()
$ws.Range("C49").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("C50"), "https://leetcode.com/problems/find-first-and-last-position-of-element-in-sorted-array/")
$ws.Range("C32").Copy()
$ws.Range("C50").PasteSpecial(-4122)

# --- Leave the selection where the author left it --------------------------
$ws.Range("G53").Select() | Out-Null
